$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 73943.664
$ws.Range("I6").Value = 39.666668
$ws.Range("J6").Value = 147847.67
$ws.Range("K6").Value = 119.000004
$ws.Range("L6").Value = 443543.01
$ws.Range("M6").Value = -7.000004000000004
$ws.Range("N6").Value = -443767.01
$ws.Range("H28").Value = 904.8
$ws.Range("I28").Value = 851.4
$ws.Range("J28").Value = 1118.4
$ws.Range("K28").Value = 851.4
$ws.Range("L28").Value = 1118.4
$ws.Range("M28").Value = -366.4
$ws.Range("N28").Value = -2088.4
$ws.Range("H33").Value = 229.91176
$ws.Range("I33").Value = 191.03703
$ws.Range("J33").Value = 379.85715
$ws.Range("K33").Value = 191.03703
$ws.Range("L33").Value = 379.85715
$ws.Range("M33").Value = 37.96297000000001
$ws.Range("N33").Value = -837.85715
$ws.Range("H64").Value = 6866.6665
$ws.Range("J64").Value = 6866.6665
$ws.Range("L64").Value = 6866.6665
$ws.Range("N64").Value = -7362.6665
$ws.Range("H67").Value = 6866.6665
$ws.Range("J67").Value = 6866.6665
$ws.Range("L67").Value = 6866.6665
$ws.Range("N67").Value = -8582.666499999999
$ws.Range("H70").Value = 797.7037
$ws.Range("I70").Value = 635.8182
$ws.Range("K70").Value = 1907.4546
$ws.Range("M70").Value = -1637.4546
$ws.Range("H73").Value = 797.7037
$ws.Range("I73").Value = 635.8182
$ws.Range("K73").Value = 1907.4546
$ws.Range("M73").Value = -971.4546
$ws.Range("H112").Value = 55557960
$ws.Range("J112").Value = 2702
$ws.Range("L112").Value = 8106
$ws.Range("N112").Value = -10322
$ws.Range("H116").Value = 1965.8334
$ws.Range("I116").Value = 1798.3334
$ws.Range("J116").Value = 2133.3333
$ws.Range("K116").Value = 1798.3334
$ws.Range("L116").Value = 2133.3333
$ws.Range("M116").Value = 1643.6666
$ws.Range("N116").Value = -9017.3333
$ws.Range("H125").Value = 1246
$ws.Range("J125").Value = 1900
$ws.Range("L125").Value = 17100
$ws.Range("N125").Value = -22020
$ws.Range("H138").Value = 3848625
$ws.Range("I138").Value = 1651
$ws.Range("J138").Value = 5885258.5
$ws.Range("K138").Value = 4953
$ws.Range("L138").Value = 17655775.5
$ws.Range("M138").Value = 187
$ws.Range("N138").Value = -17666055.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1100.174
$ws.Range("I45").Value = 1011.3333
$ws.Range("J45").Value = 1420
$ws.Range("K45").Value = 1011.3333
$ws.Range("L45").Value = 1420
$ws.Range("M45").Value = -634.3333
$ws.Range("N45").Value = -2174
$ws.Range("H61").Value = 41751376
$ws.Range("I61").Value = 58883836
$ws.Range("J61").Value = 143974.14
$ws.Range("K61").Value = 58883836
$ws.Range("L61").Value = 143974.14
$ws.Range("M61").Value = -58883624
$ws.Range("N61").Value = -144398.14
$ws.Range("H63").Value = 4659.7
$ws.Range("I63").Value = 4799.6665
$ws.Range("K63").Value = 4799.6665
$ws.Range("M63").Value = -4113.6665
$ws.Range("H66").Value = 4659.7
$ws.Range("I66").Value = 4799.6665
$ws.Range("K66").Value = 23998.3325
$ws.Range("M66").Value = -20566.3325
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H136").Value = 41751376
$ws.Range("I136").Value = 58883836
$ws.Range("J136").Value = 143974.14
$ws.Range("K136").Value = 176651508
$ws.Range("L136").Value = 431922.42
$ws.Range("M136").Value = -176648958
$ws.Range("N136").Value = -437022.42
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 1000
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -1336
$ws.Range("H25").Value = 2050.6667
$ws.Range("I25").Value = 150
$ws.Range("K25").Value = 150
$ws.Range("M25").Value = 85
$ws.Range("H105").Value = 26318522
$ws.Range("I105").Value = 45456580
$ws.Range("J105").Value = 3691.375
$ws.Range("K105").Value = 45456580
$ws.Range("L105").Value = 3691.375
$ws.Range("M105").Value = -45454833
$ws.Range("N105").Value = -7185.375
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3985.7144
$ws.Range("I62").Value = 3950
$ws.Range("K62").Value = 3950
$ws.Range("M62").Value = -3326
$ws.Range("H65").Value = 3985.7144
$ws.Range("I65").Value = 3950
$ws.Range("K65").Value = 19750
$ws.Range("M65").Value = -16630
$ws.Range("H122").Value = 2118.4814
$ws.Range("I122").Value = 1639.1428
$ws.Range("J122").Value = 2634.6924
$ws.Range("K122").Value = 4917.428400000001
$ws.Range("L122").Value = 7904.0772
$ws.Range("M122").Value = -2467.428400000001
$ws.Range("N122").Value = -12804.0772
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 2189
$ws.Range("I13").Value = 163.66667
$ws.Range("J13").Value = 3057
$ws.Range("K13").Value = 491.00001
$ws.Range("L13").Value = 9171
$ws.Range("M13").Value = -323.00001
$ws.Range("N13").Value = -9507
$ws.Range("H47").Value = 1779.6
$ws.Range("I47").Value = 559.2
$ws.Range("J47").Value = 3000
$ws.Range("K47").Value = 1677.6
$ws.Range("L47").Value = 9000
$ws.Range("M47").Value = -1246.6
$ws.Range("N47").Value = -9862
$ws.Range("H55").Value = 1458.75
$ws.Range("J55").Value = 2652.5
$ws.Range("L55").Value = 7957.5
$ws.Range("N55").Value = -8311.5
$ws.Range("H80").Value = 1401.2572
$ws.Range("I80").Value = 1061.2941
$ws.Range("J80").Value = 1722.3334
$ws.Range("K80").Value = 3183.8823
$ws.Range("L80").Value = 5167.0002
$ws.Range("M80").Value = -2247.8823
$ws.Range("N80").Value = -7039.0002
$ws.Range("H83").Value = 1401.2572
$ws.Range("I83").Value = 1061.2941
$ws.Range("J83").Value = 1722.3334
$ws.Range("K83").Value = 9551.6469
$ws.Range("L83").Value = 15501.0006
$ws.Range("M83").Value = -4871.6469
$ws.Range("N83").Value = -24861.0006
$ws.Range("H113").Value = 663.2381
$ws.Range("I113").Value = 592.1429000000001
$ws.Range("J113").Value = 805.4286
$ws.Range("K113").Value = 1776.4287
$ws.Range("L113").Value = 2416.2858
$ws.Range("M113").Value = 393.5712999999998
$ws.Range("N113").Value = -6756.2858
$ws.Range("H118").Value = 2980.4211
$ws.Range("I118").Value = 578.625
$ws.Range("J118").Value = 4727.1816
$ws.Range("K118").Value = 1735.875
$ws.Range("L118").Value = 14181.5448
$ws.Range("M118").Value = -492.875
$ws.Range("N118").Value = -16667.5448
$ws.Range("H122").Value = 880.13513
$ws.Range("J122").Value = 1294.826
$ws.Range("L122").Value = 11653.434
$ws.Range("N122").Value = -16553.434
$ws.Range("H131").Value = 902.4211
$ws.Range("J131").Value = 932.8302
$ws.Range("L131").Value = 2798.4906
$ws.Range("N131").Value = -12878.4906
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 37348.84
$ws.Range("I70").Value = 54925.75
$ws.Range("J70").Value = 5390.8184
$ws.Range("K70").Value = 54925.75
$ws.Range("L70").Value = 5390.8184
$ws.Range("M70").Value = -54655.75
$ws.Range("N70").Value = -5930.8184
$ws.Range("H73").Value = 37348.84
$ws.Range("I73").Value = 54925.75
$ws.Range("J73").Value = 5390.8184
$ws.Range("K73").Value = 54925.75
$ws.Range("L73").Value = 5390.8184
$ws.Range("M73").Value = -53989.75
$ws.Range("N73").Value = -7262.8184
$ws.Range("H80").Value = 3827.1052
$ws.Range("I80").Value = 2992.5
$ws.Range("J80").Value = 3925.2942
$ws.Range("K80").Value = 2992.5
$ws.Range("L80").Value = 3925.2942
$ws.Range("M80").Value = -1994.5
$ws.Range("N80").Value = -5921.2942
$ws.Range("H83").Value = 3827.1052
$ws.Range("I83").Value = 2992.5
$ws.Range("J83").Value = 3925.2942
$ws.Range("K83").Value = 14962.5
$ws.Range("L83").Value = 19626.471
$ws.Range("M83").Value = -9970.5
$ws.Range("N83").Value = -29610.471
$ws.Range("H123").Value = 22248.2
$ws.Range("J123").Value = 22248.2
$ws.Range("L123").Value = 22248.2
$ws.Range("N123").Value = -27148.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 31252350
$ws.Range("I7").Value = 41668884
$ws.Range("J7").Value = 2752.5
$ws.Range("K7").Value = 41668884
$ws.Range("L7").Value = 2752.5
$ws.Range("M7").Value = -41668772
$ws.Range("N7").Value = -2976.5
$ws.Range("H22").Value = 749.42426
$ws.Range("I22").Value = 715.2778
$ws.Range("K22").Value = 715.2778
$ws.Range("M22").Value = -420.2778
$ws.Range("H27").Value = 749.42426
$ws.Range("I27").Value = 715.2778
$ws.Range("K27").Value = 715.2778
$ws.Range("M27").Value = -608.2778
$ws.Range("H61").Value = 2618.76
$ws.Range("I61").Value = 2698.5557
$ws.Range("J61").Value = 2413.5715
$ws.Range("K61").Value = 2698.5557
$ws.Range("L61").Value = 2413.5715
$ws.Range("M61").Value = -2496.5557
$ws.Range("N61").Value = -2817.5715
$ws.Range("H93").Value = 1265
$ws.Range("I93").Value = 1297.6471
$ws.Range("K93").Value = 1297.6471
$ws.Range("M93").Value = -49.64709999999991
$ws.Range("H100").Value = 1800.1765
$ws.Range("I100").Value = 1320.6
$ws.Range("K100").Value = 1320.6
$ws.Range("M100").Value = -779.5999999999999
$ws.Range("H113").Value = 2618.76
$ws.Range("I113").Value = 2698.5557
$ws.Range("J113").Value = 2413.5715
$ws.Range("K113").Value = 2698.5557
$ws.Range("L113").Value = 2413.5715
$ws.Range("M113").Value = -528.5556999999999
$ws.Range("N113").Value = -6753.5715
$ws.Range("H126").Value = 31252350
$ws.Range("I126").Value = 41668884
$ws.Range("J126").Value = 2752.5
$ws.Range("K126").Value = 125006652
$ws.Range("L126").Value = 8257.5
$ws.Range("M126").Value = -125004182
$ws.Range("N126").Value = -13197.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H107").Value = 427.24243
$ws.Range("I107").Value = 333.5
$ws.Range("J107").Value = 614.7273
$ws.Range("K107").Value = 1000.5
$ws.Range("L107").Value = 1844.1819
$ws.Range("M107").Value = 919.5
$ws.Range("N107").Value = -5684.1819
